# Change jpg capitalization: "stimuli/blank.JPG" -> "stimuli/blank.jpg"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$used = $ws.UsedRange
$found = $used.Find("stimuli/blank.JPG", [Type]::Missing, -4123, 1)

if ($found -ne $null) {
    $firstAddress = $found.Address()
    while ($true) {
        $found.Value2 = "stimuli/blank.jpg"
        $found = $used.FindNext($found)
        if ($found -eq $null -or $found.Address() -eq $firstAddress) {
            break
        }
    }
}

$wb.Save()
